# Generate Report for Handoff
# Adds a new tracked file "e8088aa3-5fd2-41a5-8060-ea3b75c18b96.md" (status:
# "Ready for handoff", not yet handed back) as a new row at the bottom of
# each of the three worksheets (Overview, zh-cn, de-de), mirroring the
# layout already used for the other "Ready for handoff" rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A9").Value = "e8088aa3-5fd2-41a5-8060-ea3b75c18b96.md"
$ws1.Range("B9").Value = "Ready for handoff"
$ws1.Range("C9").Value = "Ready for handoff"
$ws1.Range("D9").Value = "2016-03-22 00:38:57"

$ws1.Hyperlinks.Add($ws1.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/e8088aa3f5d241a58060ea3b75c18b96e8088aa3/e2e/e8088aa3-5fd2-41a5-8060-ea3b75c18b96.md", "", "", "e8088aa3-5fd2-41a5-8060-ea3b75c18b96.md")

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A9").Value = "e8088aa3-5fd2-41a5-8060-ea3b75c18b96.md"
$ws2.Range("B9").Value = ".md"
$ws2.Range("C9").Value = "Ready for handoff"
$ws2.Range("D9").Value = "e8088aa3-5fd2-41a5-8060-ea3b75c18b96.e1568a30eeff22474690f39448d8625f901cd9e9.zh-cn.xlf"
$ws2.Range("E9").Value = "2016-03-22 00:38:54"
$ws2.Range("H9").Value = "0001-01-01 00:00:00"
$ws2.Range("J9").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/e8088aa3f5d241a58060ea3b75c18b96e8088aa3/e2e/e8088aa3-5fd2-41a5-8060-ea3b75c18b96.md", "", "", "e8088aa3-5fd2-41a5-8060-ea3b75c18b96.md")
$ws2.Hyperlinks.Add($ws2.Range("D9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e8088aa3f5d241a58060ea3b75c18b96e8088aa3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e8088aa3-5fd2-41a5-8060-ea3b75c18b96.e1568a30eeff22474690f39448d8625f901cd9e9.zh-cn.xlf", "", "", "e8088aa3-5fd2-41a5-8060-ea3b75c18b96.e1568a30eeff22474690f39448d8625f901cd9e9.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A9").Value = "e8088aa3-5fd2-41a5-8060-ea3b75c18b96.md"
$ws3.Range("B9").Value = ".md"
$ws3.Range("C9").Value = "Ready for handoff"
$ws3.Range("D9").Value = "e8088aa3-5fd2-41a5-8060-ea3b75c18b96.e1568a30eeff22474690f39448d8625f901cd9e9.de-de.xlf"
$ws3.Range("E9").Value = "2016-03-22 00:38:57"
$ws3.Range("H9").Value = "0001-01-01 00:00:00"
$ws3.Range("J9").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/e8088aa3f5d241a58060ea3b75c18b96e8088aa3/e2e/e8088aa3-5fd2-41a5-8060-ea3b75c18b96.md", "", "", "e8088aa3-5fd2-41a5-8060-ea3b75c18b96.md")
$ws3.Hyperlinks.Add($ws3.Range("D9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e8088aa3f5d241a58060ea3b75c18b96e8088aa3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e8088aa3-5fd2-41a5-8060-ea3b75c18b96.e1568a30eeff22474690f39448d8625f901cd9e9.de-de.xlf", "", "", "e8088aa3-5fd2-41a5-8060-ea3b75c18b96.e1568a30eeff22474690f39448d8625f901cd9e9.de-de.xlf")

Write-Output "Report row appended for e8088aa3-5fd2-41a5-8060-ea3b75c18b96 across Overview/zh-cn/de-de"
